# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column (G) values calculated/re-derived for each data row (2-27)
$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 4
    7  = 2
    8  = 3
    9  = 2
    10 = 4
    11 = 1
    12 = 0
    13 = 1
    14 = 1
    15 = 5
    16 = 3
    17 = 0
    18 = 0
    19 = 4
    20 = 1
    21 = 4
    22 = 8
    23 = 1
    24 = 3
    25 = 5
    26 = 3
    27 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
